# Adding, editing, deleting Products and Artists
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Product "Test 1"): Artist ID 8 -> 11, Genre ID 35 -> 1
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 1

# Row 3 (Product "Test2"): Artist ID 9 -> 12 (Genre ID stays 2)
$ws.Range("A3").Value = 12

# Row 4 (Product "Test 3"): Artist ID 4324 -> 14, Genre ID 77 -> 2
$ws.Range("A4").Value = 14
$ws.Range("B4").Value = 2

# Selection ends up on the whole of row 1 (header row)
$ws.Range("A1:XFD1").Select()
